$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($addr, $value) {
    # Force the cell to remain text (not get auto-converted to a number)
    # without leaving a visible number-format style behind.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "59.790.60"
$ws.Range("E2").Value = "  +2.57%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.414.76"
$ws.Range("E3").Value = "  +1.91%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
Set-Text "D5" "550.21"
$ws.Range("E5").Value = "  +0.77%  "

# Row 6 - Solana
Set-Text "D6" "136.85"
$ws.Range("E6").Value = "  +2.26%  "

# Row 7 - USDC
Set-Text "D7" "1.00"
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - XRP
Set-Text "D8" "0.584"
$ws.Range("E8").Value = "  +3.03%  "

# Row 9 - Dogecoin
Set-Text "D9" "0.105"
$ws.Range("E9").Value = "  -0.83%  "

# Row 10 - Toncoin
Set-Text "D10" "5.72"
$ws.Range("E10").Value = "  +3.16%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -2.07%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -0.32%  "

# Row 13 - Avalanche
Set-Text "D13" "24.71"
$ws.Range("E13").Value = "  +2.75%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.845.01"
$ws.Range("E14").Value = "  +1.93%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "59.770.07"
$ws.Range("E15").Value = "  +2.65%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +0.06%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.433.37"
$ws.Range("E17").Value = "  +1.42%  "

# Row 18 - Chainlink
Set-Text "D18" "11.27"
$ws.Range("E18").Value = "  +2.54%  "

# Row 19 - Polkadot
Set-Text "D19" "4.36"
$ws.Range("E19").Value = "  +0.55%  "

# Row 20 - BitcoinCash
Set-Text "D20" "329.98"
$ws.Range("E20").Value = "  -0.45%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -2.93%  "

# Row 22 - Dai
Set-Text "D22" "0.998"

# Row 23 - Litecoin
Set-Text "D23" "65.78"
$ws.Range("E23").Value = "  +3.63%  "

# Row 24 - Kaspa
$ws.Range("E24").Value = "  +2.51%  "

# Row 25 - InternetComputer(DFINITY)
Set-Text "D25" "8.57"
$ws.Range("E25").Value = "  +3.49%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("E26").Value = "  +0.16%  "

# Row 27 - Fetch.AI
$ws.Range("E27").Value = "  +0.91%  "

# Row 28 - PEPE
$ws.Range("D28").Value = "0.0₃0775"
$ws.Range("E28").Value = "  +3.94%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  +0.53%  "

# Row 30 - Monero
Set-Text "D30" "170.56"
$ws.Range("E30").Value = "  +0.02%  "

# Row 31 - Aptos
Set-Text "D31" "6.17"
$ws.Range("E31").Value = "  +0.05%  "

# Row 32 - EthereumClassic
$ws.Range("E32").Value = "  +0.81%  "

# Row 33 - SuiNetwork
$ws.Range("E33").Value = "  +1.38%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  +0.16%  "

# Row 37 - NEARProtocol
Set-Text "D37" "4.16"
$ws.Range("E37").Value = "  -0.71%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  -0.71%  "

# Row 39 - OKB
Set-Text "D39" "39.39"
$ws.Range("E39").Value = "  +0.57%  "

# Rows 40/41 swap: Bittensor moves to row 40, PolygonEcosystemToken moves to row 41
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-Text "D40" "315.14"
$ws.Range("E40").Value = "  +9.42%  "

$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-Text "D41" "0.409"
$ws.Range("E41").Value = "  -1.74%  "

# Row 42 - Filecoin
$ws.Range("E42").Value = "  -1.03%  "

# Row 43 - Aave
Set-Text "D43" "137.83"
$ws.Range("E43").Value = "  -3.22%  "

# Row 44 - Stellar
Set-Text "D44" "0.0962"
$ws.Range("E44").Value = "  +1.09%  "

# Row 45 - Hedera
Set-Text "D45" "0.0517"
$ws.Range("E45").Value = "  -0.51%  "

# Rows 46/47 swap: InjectiveProtocol moves to row 46, Mantle moves to row 47
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-Text "D46" "19.35"
$ws.Range("E46").Value = "  +1.85%  "

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-Text "D47" "0.578"
$ws.Range("E47").Value = "  +2.03%  "

# Row 48 - Polygon
Set-Text "D48" "0.406"
$ws.Range("E48").Value = "  +4.83%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  +0.26%  "

# Row 50 - EnergySwap
Set-Text "D50" "17.52"
$ws.Range("E50").Value = "  -0.28%  "

# Row 51 - WhiteBITCoin
$ws.Range("E51").Value = "  -0.35%  "
